# Fruta / hortaliza, semanal
# Inserts a new week's worth of data (4 rows) for
# "Terminal Hortofrutícola Agro Chillán - Naranja" right before the
# previous most-recent week (row 215), pushing the existing rows down by
# four. The new rows reuse the same product/quality/price metadata as the
# (old) most-recent week, but carry the new reporting date and volumes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the current first data block (rows 215-218),
# which shifts the old rows 215-260 down to 219-264 (Excel-style insert:
# existing formatting/styles on column D carry over automatically).
$ws.Rows("215:218").Insert()

# New week's date (was 44468, now 44476 == 2021-10-07).
$newDate = 44476

$ws.Range("A215").Value = 7
$ws.Range("B215").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C215").Value = "Ñuble"
$ws.Range("D215").Value = $newDate
$ws.Range("E215").Value = 16
$ws.Range("F215").Value = "Fruta"
$ws.Range("G215").Value = 100102
$ws.Range("H215").Value = "Cítricos"
$ws.Range("I215").Value = 100102005
$ws.Range("J215").Value = "Naranja"
$ws.Range("K215").Value = "Lane Late"
$ws.Range("L215").Value = "Primera"
$ws.Range("M215").Value = 300
$ws.Range("N215").Value = 6000
$ws.Range("O215").Value = 6500
$ws.Range("P215").Value = 6250
$ws.Range("Q215").Value = "$/bandeja 15 kilos granel"
$ws.Range("R215").Value = "Región de O'Higgins"
$ws.Range("S215").Value = 417
$ws.Range("T215").Value = 15

$ws.Range("A216").Value = 7
$ws.Range("B216").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C216").Value = "Ñuble"
$ws.Range("D216").Value = $newDate
$ws.Range("E216").Value = 16
$ws.Range("F216").Value = "Fruta"
$ws.Range("G216").Value = 100102
$ws.Range("H216").Value = "Cítricos"
$ws.Range("I216").Value = 100102005
$ws.Range("J216").Value = "Naranja"
$ws.Range("K216").Value = "Lane Late"
$ws.Range("L216").Value = "Segunda"
$ws.Range("M216").Value = 240
$ws.Range("N216").Value = 5000
$ws.Range("O216").Value = 5500
$ws.Range("P216").Value = 5250
$ws.Range("Q216").Value = "$/bandeja 15 kilos granel"
$ws.Range("R216").Value = "Región de O'Higgins"
$ws.Range("S216").Value = 350
$ws.Range("T216").Value = 15

$ws.Range("A217").Value = 7
$ws.Range("B217").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C217").Value = "Ñuble"
$ws.Range("D217").Value = $newDate
$ws.Range("E217").Value = 16
$ws.Range("F217").Value = "Fruta"
$ws.Range("G217").Value = 100102
$ws.Range("H217").Value = "Cítricos"
$ws.Range("I217").Value = 100102005
$ws.Range("J217").Value = "Naranja"
$ws.Range("K217").Value = "Navel Late"
$ws.Range("L217").Value = "Primera"
$ws.Range("M217").Value = 240
$ws.Range("N217").Value = 6000
$ws.Range("O217").Value = 6500
$ws.Range("P217").Value = 6250
$ws.Range("Q217").Value = "$/bandeja 15 kilos granel"
$ws.Range("R217").Value = "Región de O'Higgins"
$ws.Range("S217").Value = 417
$ws.Range("T217").Value = 15

$ws.Range("A218").Value = 7
$ws.Range("B218").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C218").Value = "Ñuble"
$ws.Range("D218").Value = $newDate
$ws.Range("E218").Value = 16
$ws.Range("F218").Value = "Fruta"
$ws.Range("G218").Value = 100102
$ws.Range("H218").Value = "Cítricos"
$ws.Range("I218").Value = 100102005
$ws.Range("J218").Value = "Naranja"
$ws.Range("K218").Value = "Navel Late"
$ws.Range("L218").Value = "Segunda"
$ws.Range("M218").Value = 240
$ws.Range("N218").Value = 5000
$ws.Range("O218").Value = 5500
$ws.Range("P218").Value = 5250
$ws.Range("Q218").Value = "$/bandeja 15 kilos granel"
$ws.Range("R218").Value = "Región de O'Higgins"
$ws.Range("S218").Value = 350
$ws.Range("T218").Value = 15
